$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values (Part Code / Unit text changed, Rate/Qty numbers changed) ---
$ws.Range("A2").Value = "0001-908"
$ws.Range("B2").Value = "aqa"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# --- New columns populated on row 2: Ori. MRP (text), PO No. (number), Date (date) ---
$ws.Range("E2").Value = "dsf"
$ws.Range("F2").Value = 1

$ws.Range("G2").NumberFormat = "mm-dd-yy"
$dateValue = Get-Date -Year 2024 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("G2").Value = $dateValue

# --- Column G width to fit the date column ---
$ws.Columns.Item(7).ColumnWidth = 9.6666666667

# --- Update the active selection to match the saved view state ---
$ws.Range("G12").Select()

$wb.Save()
